$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'Bang bros, b-bang bros Bang bros, b-bang bros Bang bros, b-bang, bros Bang bros, b-bang bros 3, 2, 1... Joke, yeah Je n''ai pas de temps à perdre, négro, quand j''passe au micro j''ai mauvais caractère Tout est parallèle dans mes structures, je ne calcule pas la Terre, j''suis le futur Négro j''suis plus dur qu''un pare choc de tu-ture, yeah faut l''cutter au point comme une suture Tout ou rien, rien, les rappeurs au pied comme des toutous, hein, hein Les MC sont dingues, leurs meufs sont dingues de moi Je suis dur, très dur, encore plus dur que tes fins de mois Les gars qui m''comprennent pas pensent que j''écris d''la merde OK, crise d''épilepsie quand j''les éclaire, ils ont les nerfs OK, yeah, yeah, yeah, yeah ! Disons qu''j''ai l''aura que t''auras pas, négro ouvre les yeux tu me rateras pas Le jour où les poules iront au KFC j''suis sûr, j''suis sûr tu me rattraperas Les meufs ne nous font pas trop chier Par contre les bombes deviennent le contraire de Castro, yeah Me prends pas trop pour un nigaud avec ta fausse paire de Bape petit J''suis pas Pharrell, dis à ta mère qu''elle t''achète mon sque-di J''suis pas Kid Cudi, j''sais pas qui t''l''a dit mais quitte-là Si ta meuf écoute mon ''sque-d c''est'' qu''elle te trompe, mec quitte-la Quand on parle, c''est Bang Bros, quand on rappe, c''est Bang Bros Quand on sape, c''est Bang Bros Bang, Bang, Bang, Bang Bros Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Quand on parle, c''est Bang Bros, quand on rappe, c''est Bang Bros Quand on sape, c''est Bang Bros Bang, Bang, Bang, Bang Bros Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Les rappeurs s''branlent sur leur style mais en vrai, ils lèchent les vitrines Moi tu sais qu''c''est Gre-Gre-Greane, je ne lècherai qu''ta cyprine À chaque rime, ça fait une victime Écoute, j''fais beaucoup d''rimes donc ça fait beaucoup de victimes Allez, fait tourner la weed, fais nunu la ''teille de whisky Ta meuf est bonne, elle a les yeux bleus, le poil doux comme un husky J''ai vu sous sa culotte, elle a les bouclettes à la Starsky Dis-moi négro, n''y aurait-il pas erreur de casting C''est pour tous mes mecs à court de mots-clés dans Youporn Pour tous mes mecs en sang sur J.Lo dans U-Turn J''kifferai son tarma jusqu''à mon enterrement Des meufs, du flow, d''la thune, tout ça, tu sais qu''j''en ai tellement Ok, c''est Joke et ton pote, Gre-Gre-Greane Juste un aller-retour sur le gri-gri-grill J''renvoies tous ses branleurs au bercail Au bercail Ok, t''es trop sympa'' mais t''as le swag à Orelsan Quand on parle, c''est Bang Bros, quand on rappe, c''est Bang Bros Quand on sape, c''est Bang Bros, Bang, Bang, Bang, Bang Bros Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Quand on parle, c''est Bang Bros, quand on rappe, c''est Bang Bros Quand on sape, c''est Bang Bros, Bang, Bang, Bang, Bang Bros Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Yeah, lève ton verre en l''air, lè-lè-lève ton verre en l''air Pour tous mes mecs en sang sur J-Lo dans U-Turn Je kifferais son terma, jusqu''à mon enterrement Des meufs, du flow, de la thune, tout ça tu sais que j''en ai tellement'

$ws.Range("C3").Value = 'Yeah! OK! Charly Greane, Bx! James Heartbreaker Bordeaux! BX! OK Bx, c''est la ville la référence, d''évidence ta préférence yeah J''plante mon drapeau sur la carte, yeah Toutes les rues, tous les quartiers Ok mais mec ici, on rappe comme on parle Et quand on parle trop on t''chasse Sale Sud avec des sales têtes des sales putes des sales mecs des sales rues des salles pleines Dans ma ville, on te malmène OK mais Greane, on m''a dit, tu te la pètes on m''a dit tu fais pas l''poids On m''a médit mais dans sa aah j''avais déjà mes 10 doigts 2 OK tu l''a trouvée, ma ville tu la voudrais Fais bander les bandits B.X Quand on arrive dans ta ville, toutes tes amies veulent nous serrer Fais bander les bandits B.X Mec ici c''est sea sex and sun Barbecue, Greane, c''est ce qu''elles veulent OK on est Frais dans nos tires, on est Prêts à sortir et la nuit, on ment comme on respire Laisse faire les experts, laisse faire ils espèrent mais sont forcés à se taire On est trop cotés, on t''a trop choqué Tu voudrais boycotter mais ta copine a ma photo sous l''oreiller On a trop d''admirateurs Les filles mettent ma musique pendant qu''elles passent l''aspirateur 2 2 OK, dans ma ville, on t''aime pas ça baise comme dans Ken Park ma ville c''est l''extase On fait feu sur toi comme sur un gardien de handball Quand on débarque sur la place avec une grosse paire de Nike! 21'

$ws.Range("C4").Value = 'On est Tes amis Tu n''te sentiras plus jamais seul, vas-y Il est 6 heures 30, il est 6 heures 30 Faut aller en cours, ton réveille sonne ton réveille sonne Ta mère grogne, et sa voix résonne sa voix résonne Mais que justice soit faite, matte comment l''école te traite Il est 8 heures du mat'' et t''es dehors en plein hiver pour un cours d''EPS cours d''EPS Gèle-toi les fesses gèle-toi les fesses Mais quel que soit la saison Multiples sont les raisons Tu t''es couché trop tard ouais tu t''es couché trop tard Ou tu n''as pas fais tes d''voirs tu n''as pas fait tes d''voirs Reproduis la signature d''tes parents Continue oué j''te l''assure c''est marrant Et si les profs tiennent ta meuf prisonnière Demande-lui aussi de faire l''école buissonnière x2 Les cours on s''en fout, on les rate on les rate Pour aller fourrer des Kicker des oh C''est pas l''école qui va me trouver du taf non T''as cru qu''j''allais passer ma vie à bouffer des pâtes non Greane, Tekitek ohh On est on est Tes amis oui Tu n''te sentiras plus jamais seul, vas-y Manquer un devoir ne pas se faire culpabiliser culpabiliser Sécher ta journée sera rentabilisée rentabilisée Marre-toi maintenant, tu pleureras plus tard Tu s''ras jamais seul dans ç''cas-là Y''en a un qui sèche tout le monde suit c''est l''escalade Des friandises pleins ton cartable pleins ton cartable Et des meufs plein ton portable plein ton portable C''est l''heure de rentrer à la maison Ton père est d''vant la télévision Toi tu fais profil bas toi tu fais profil bas S''il te teste ne rougis pas s''il te teste ne rougis pas Je te jure qu''être immature c''est marrant Prendre des rides en vieillissant c''est navrant À la fin du trimestre, surveille ta boîte aux lettres Enfoiré, t''as fait l''école buissonière On est on est Tes amis oui Tu n''te sentiras plus jamais seul vas-y'

$ws.Range("C5").Value = 'Hein, Charly Greane Western Love Ok Jeunes et fiers We are much too young, and life''s so big We don''t know yet what the future brings Any sense for us on est jeunes Any sense for us on est fiers We are much too young, and life''s so big We don''t know yet what the future brings Any sense for us on est jeunes Any sense for us et on est fiers Parce-qu''on est jeunes et fiers Parce-qu''ils sont vieux et ternes 3 mots sur un papier et c''est leur musique qu''on enterre ouais Oui y''a longtemps tu as eu ton heure de gloire ton heure de gloire Il y''a longtemps, très longtemps Alors voilà la relève, la nouvelle ère, le sang neuf Et tu t''en mort les doigts, ouais il fallait nous tuer dans l''uf Tu veux suivre le train suivre le train , mais rien y fait rien ? Il y''a longtemps que le train t''a dépassé Et tu n''trompes personne, non l''habit n''fait pas l''moine Il suffit pas d''une casquette pour tromper sur ton âge sur ton âge Moi je n''réfléchis pas au lendemain au lendemain J''avance tant pis si j''perds des plumes en chemin Je suis jeune et fier, toi es-tu jeune et fier ? Nous, nous sommes jeunes et fiers mais êtes-vous jeunes et fiers ? Je-je suis jeune je suis jeune Et fier Mais toi, es-tu jeune ? Es-tu jeune et fier ? Laisse-les parler, et laisse-les ramer Mais une fois tout en haut regarde-les s''écraser La puissance est dans tes mains, ouais il faut y croire J''viens pour les récompenses comme Marion Cotillard Ma musique s''écoute au lycée, au collège, ou à la fac Toi t''es largué, tu comptes encore vendre tes CDs à la fnac C''est fini, comme ton époque c''est fini J''collectionne les médailles tu collectionnes les vinyles On est jeunes on est jeunes, et on est fiers de l''être Donne-moi une paire de lèvres Mauvais élève vient à ta table et prend l''pouvoir Un bruit d''couloir, invite ta team dans un mouloir Je suis jeune et fier toi es-tu jeune et fier ? Nous, nous sommes jeunes et fiers mais êtes-vous jeunes et fiers ? Je-je suis jeune je suis jeune Et fier Mais toi, es-tu jeune ? Es-tu jeune et fier ? We are much too young, and life''s so big We don''t know yet what the future brings Any sense for us on est jeunes Any sense for us on est fiers We are much too young, and life''s so big We don''t know yet what the future brings Any sense for us Young and proud, we are'

$ws.Range("C6").Value = 'Mesdames, mesdemoiselles, messieurs S''il vous plaît Faites du bruit pour Orgasmic Balance la sauce, Gaga Charly Greane J''ai 20 piges, mais jusqu''au Japon On achète mes disques Ils savent même où j''habite Alors je brouille les pistes J''gribouille trois mots sur une page et ta carrière n''existe plus oh, oh! J''ai du respect pour les femmes j''te traite quand même de fils de pute Oh, oh ! Fils de, fils de pute ! Direct de la cave, ok Tous ces bâtards m''écoutent et répètent tout comme des peroquets, ok! J''rentre dans la boîte et tous les gars s''niquent ok, ok On m''donne des thunes contre un bout d''plastique ok, ok Et toutes les meufs me trouvent trop orgasmique Fais une croix sur ton rôle si j''suis dans l''casting Habillé, j''suis l''nouveau Bertrand Cantat En casquette New Era Tous mes raps Sont rentables Tu verras Habillé, j''suis l''nouveau Bertrand Cantat En casquette New Era Tous mes raps Sont rentables Tu verras Refrain - Cuizinier On a de vraiment jolies Nike mais y''a pas qu''ça qui tue chez nous Des cerveaux gigantesques c''est Cuiz et Greane qui tuent, c''est tout ! De vraiment jolies Nike mais y''a pas qu''ça qui tue chez nous Des cerveaux gigantesques c''est Greane et Cuiz qui tuent, c''est fou ! Cuizinier 7, 5, double 3 Cuiz et Greane te donnent chaud comme une couche de gras Les femmes aiment nous caresser, nous sommes des boules de poil Si tu es jaloux, tu n''aimes pas, alors bouge de là On a pas d''temps à perdre nan ! D''l''argent à s''faire Un mec comme moi va dans l''magasin et achète 4 paires Les rythmes de nos vies sont excessifs Quand tu nous vois tes pensées deviennent dépressives C''est ton pote Cuizi Cuiz, Greane, Gaga et Felix On m''connaît partout dans l''monde ma musique est un satellite Au moins une radio dans chaque pays qui passe mes disques Je prends la pose, fais un sourire et lâche des bises On espère tous ici qu''t''apprécies cette chanson On a besoin d''ton fric pour payer à nos femmes des pensions Charly Greane, Cuizi Cuiz, 2007 On sait plus où donner d''la tête comme si on avait plus d''2 000 sexes Refrain - Charly Greane Charly Greane Ok c''est Charly Greane et Cuizi et Cuizi On renvoie tous les MCs défectueux à l''usine à l''usine Mais matte ces mecs qui croient trop qui sont dans l''coup coup J''les grille direct comme un suçon dans l''cou cou Cuizinier Coucou, c''est Cuiz, comment ça va ma cocotte ? Tu veux v''nir dans ma cuisine et savoir ce que j''y concocte ? Tu gout''ras p''têt'' à ma prochaine recette dans Brain On fait plusieurs trafics tout ça sous la même enseigne Charly Greane Mec mais en soirée, viens pas gratter l''amitié On fait d''la musique pour les stades comme Johnny Hallyday Quand on rentre dans le club, c''est pour boire et s''amuser Charly Greane arrive trop frais, laissez entrer l''accusé Cuizinier On m''accuse de plein d''choses mais c''est d''la connerie en boîte On m''reproche d''aimer les filles, les sapes et d''sortir en boîte Je suis désolé nous n''avons pas les mêmes valeurs Jamais en r''tard, trop en avance, et j''essaye d''être à l''heure Refrain 2x'

$ws.Range("C7").Value = 'Pourquoi faudrait-il choisir une fille lorsque l''on peut en avoir une Je n''arrive pas à me décider Je n''ai pas envie de te blesser Je m''excuse aux près de celles à qui j''aurais pu faire du mal Je suis désolé mais mon cur est brisés ok C''est pour cela que j''aime quand tu bouges comme sa C''est pour cela que j''aime quand tu remues ton cul Et danse danse sur la piste et sur ma bite C''est Yung Sid et Charly Greane L''orgasme est assuré x2 Blow on me Blow on Blow x2 Touch x2 x8 Blow on me Touch me Touch x8 Touch me x2 Caress me Touch me Caress me Touche moi,caresse moi et aime moi Déteste moi mais suce la Suce moi C''est Charly Oui c''est Charly qu''elles acclament Mon pénis est dans ta femme Et tout mes doigts sont dans sa chatte C''est Charly le roi du charme Toutes les putes se noient dans mon sperme Sa bouche s''ouvre et se referme Ses trous s''ouvrent quand je la baise C''est brutal tu te dis merde Mais touche moi Caresse moi Tou-Touche moi Caresse moi donc touche moi Caresse moi Oui pétasse soit éfficace x2 Blow on me Blow on Blow x2 Touch x2 x8 Blow on me Touch me Touch x8 Touch me x2 Caress me Touch me Caress me Bite bite fesses C''est vite vite fais ça Ton fessier bounce Ne fait pas semblant Tout les bonhommes bandent Charly G. est dans la salle avec son nombril met l''ambiance Tout les gros culs en évidence C''est sur Du sperme sur ton visage Ta l''air beaucoup moins sage Mais t''inquiète c''est n''est qu''un passage Plus tard on passe aux massages Et c''est ma bite bite dans tes fesses fesses qui fais de toi une princesse J''aime quand tu rebondis sur le tapis x3 Blow on me Blow on Blow x2 Touch x2 x8 Blow on me Touch me Touch x8 Touch me x2 Caress me Touch me Caress me'
